# MAJ du texte SYSML A faire
#
# The diagram of connectors around the "Compte tours" sensor (top-right
# cluster, slide 6) is resized/repositioned: the horizontal dashed leader
# lines get shorter and the small bracket/arrow cluster moves left & up
# slightly. One connector ("Connecteur droit 10") also changes its
# stacking position (it now sits right after "Connecteur droit 11"
# instead of right before it).
#
# PowerPoint's COM object model only exposes shape geometry in points
# (Shape.Left / .Top / .Width / .Height) while the underlying OOXML keeps
# EMU (1 pt = 12700 EMU); the point values below were chosen (each is an
# exact IEEE-754 float32 value) so that converting them back to EMU
# reproduces the exact target integer EMU from the reference file.

function Set-ShapeGeometry($shape, $left, $top, $width, $height) {
    $shape.Left   = $left
    $shape.Top    = $top
    $shape.Width  = $width
    $shape.Height = $height
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# "Connecteur droit 7" (id 8) - dashed leader from the rounded-rect label;
# off unchanged (5652120,3425791), ext cx 3240360 -> 2520360.
Set-ShapeGeometry $s.Shapes.Item(1) 445.048828125 269.7473449707031 198.45355224609375 0.0

# "Connecteur droit 10" (id 11) - off (7019273,1619703) -> (7024846,1681417),
# ext (1427655,0) -> (1005517,0).
Set-ShapeGeometry $s.Shapes.Item(2) 553.1375122070312 132.39505004882812 79.17456817626953 0.0

# "Connecteur droit 11" (id 12) - off unchanged (6827904,1835775),
# ext cx 2136584 -> 1344576.
Set-ShapeGeometry $s.Shapes.Item(3) 537.6302490234375 144.54922485351562 105.87213134765625 0.0

# "Connecteur droit 12" (id 13) - off (7013142,2051847) -> (7020528,1985782),
# ext (1433786,0) -> (1009835,0).
Set-ShapeGeometry $s.Shapes.Item(4) 552.7974853515625 156.3607940673828 79.51457214355469 0.0

# "Connecteur droit 13" (id 14) - off (8446928,1619703) -> (8030363,1681417),
# ext (0,432144) -> (0,304365).
Set-ShapeGeometry $s.Shapes.Item(5) 632.3120727539062 132.39505004882812 0.0 23.965749740600586

# "Connecteur droit 14" (id 15) - off (8446928,1619703) -> (8030363,1681417),
# ext (84465,216120) -> (59490,152216).
Set-ShapeGeometry $s.Shapes.Item(6) 632.3120727539062 132.39505004882812 4.6842522621154785 11.985512733459473

# "Connecteur droit 15" (id 16) - off (8446928,1835775) -> (8030363,1833600),
# ext (84465,216072) -> (59490,152182).
Set-ShapeGeometry $s.Shapes.Item(7) 632.3120727539062 144.37796020507812 4.6842522621154785 11.982834815979004

# "Connecteur droit 16" (id 17) - off (7019225,1547743) -> (7024812,1630735),
# ext (1224136,48) -> (862176,34).
Set-ShapeGeometry $s.Shapes.Item(8) 553.1348266601562 128.4043426513672 67.88787841796875 0.002677165437489748

# "Connecteur droit 17" (id 18) - off (7019177,2123807) -> (7024779,2036464),
# ext (1224184,48) -> (862210,34).
Set-ShapeGeometry $s.Shapes.Item(9) 553.1322631835938 160.35150146484375 67.89055633544922 0.002677165437489748

# "Connecteur droit 18" (id 19) - off (8243361,1547743) -> (7886988,1630735),
# ext (0,576064) -> (0,405729).
Set-ShapeGeometry $s.Shapes.Item(10) 621.022705078125 128.4043426513672 0.0 31.947166442871094

# Finally, re-stack "Connecteur droit 10" (currently shape #2) one step
# forward so it ends up right after "Connecteur droit 11" (matching the
# new shape order in the reference XML): msoBringForward = 2.
$s.Shapes.Item(2).ZOrder(2)
